$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "58.626.77"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.08%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.481.65"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("E4").Value = "  -0.21%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "529.89"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.23%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "132.31"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.02%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.69%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.564"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.93%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.480.67"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.97%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0991"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.33%  "

$ws.Range("E11").Value = "  -2.95%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.10"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.01%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.326"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.18%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.935.09"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.01%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "58.492.31"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.61%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "22.08"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000134"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.03%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.505.09"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.28%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "10.52"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.34%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.21"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.38%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "319.46"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.20%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.18"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.74%  "

$ws.Range("E23").Value = "  +0.22%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "65.45"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.80%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.404"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.21%  "

$ws.Range("E27").Value = "  -0.55%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.33"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "174.98"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.38%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0₃0746"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.38%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.71"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.04%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.17"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.27%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "6.19"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("E34").Value = "  +0.17%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.60%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "17.98"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.21"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -3.50%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.87"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.94%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.48"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.81%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "36.40"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.90%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.811"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +6.03%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.45"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.99%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "270.70"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.02%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "129.36"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +7.52%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "4.93"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.45%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.587"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.13%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0934"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.76%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0502"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.08%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0216"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.69%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "16.61"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.48%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.732.25"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "
